$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: title paragraph - merge the "dimensionality reduction" run and
# the following single-space run into one run (no visible text change).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)

$xmlTitle = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="494D928A" w14:textId="4B5F1F2E" w:rsidR="00703F7F" w:rsidRPr="00703F7F" w:rsidRDefault="00703F7F" w:rsidP="00703F7F">
<w:pPr><w:pStyle w:val="a3"/><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00703F7F"><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Interpretation and analysis of </w:t></w:r>
<w:r w:rsidR="00253233" w:rsidRPr="00253233"><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">dimensionality reduction </w:t></w:r>
<w:r w:rsidRPr="00703F7F"><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>results</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r1.InsertXML($xmlTitle)

# ---------------------------------------------------------------------------
# Change 2: final analysis paragraph - rewritten with an inserted citation
# and a few phrasing tweaks ("Alpha strand" -> "Alpha chains", "Beta
# strands" -> "Beta chains", "In addition, ..." -> "We believe the reason
# is ...").
# ---------------------------------------------------------------------------
$target = "Overall, the human TCR is more dispersed"
$pFinal = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith($target)) {
        $pFinal = $para
        break
    }
}

$r2 = $d.Range($pFinal.Range.Start, $pFinal.Range.End)

$xmlFinal = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="38714DE9" w14:textId="2899F7FD" w:rsidR="00703F7F" w:rsidRPr="00703F7F" w:rsidRDefault="00703F7F" w:rsidP="00703F7F">
<w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00703F7F"><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Overall, the human TCR is more dispersed and more difficult to form a single cluster compared to the mouse. It is possible that the high genetic diversity of the human population, the more complex immune system, and the exposure to a wide variety of pathogens in different environments have caused the human TCR to exhibit a more diverse clustering behaviour. Single Alpha chain data exhibit a large amount of overlap and form small dispersed clusters; single Beta chain data are denser with distinct cluster boundaries. This may be due to the fact that during gene rearrangement, the Alpha </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">chains </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>does not have a D region and only undergoes VJ rearrangement</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>[</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Rosati, E., Dowds, C.M., </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Liaskou</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">, E. et al. Overview of methodologies for T-cell receptor repertoire analysis. BMC </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Biotechnol</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> 17, 61 (2017). https://doi.org/10.1186/s12896-017-0379-9</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>]</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">, which has a lower antigenic specificity compared to Beta. </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">And after combining Alpha and Beta, the data points appeared to be more dispersed overall. </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>We believe the reason</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> is d</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">ifferent TCRs have different combinations of features on the Alpha and Beta </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:hint="eastAsia"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>chains</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>, resulting in a more diverse distribution, making the data points more scattered after combination.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r2.InsertXML($xmlFinal)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
